# feat: add 2022-Q1 data
#
# What happened conceptually (per the source commit):
#   - The former "总计" (totals) sheet is renamed to "2022-Q1" and its
#     content is replaced by the new quarter's per-fund holding detail.
#   - A brand new "总计" sheet is created right after it, holding the
#     same rolled-up summary table as before plus a new first data row
#     for 2022-Q1 (and the old rows shifted down by one).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# Create the new summary ("总计") sheet right after the current one so the
# final tab order is: ... 2021-Q4, 2022-Q1, 总计
$newTotal = $wb.Worksheets.Add($null, $total)
$newTotal.Name = "总计_tmp_rename"

# Repurpose the existing "总计" sheet as the "2022-Q1" detail sheet, then
# name the freshly added sheet "总计".
$total.Name = "2022-Q1"
$newTotal.Name = "总计"

$detail = $wb.Worksheets.Item("2022-Q1")
$summary = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Clear the old totals content from what is now the "2022-Q1" sheet and
#    populate it with the per-fund holding detail rows.
# ---------------------------------------------------------------------
$detail.Cells.Clear()

$detail.Range("B1").Value = "基金代码"
$detail.Range("C1").Value = "基金名称"
$detail.Range("D1").Value = "基金规模"
$detail.Range("E1").Value = "股票总仓位"
$detail.Range("F1").Value = "仓位占比"
$detail.Range("G1").Value = "持有市值(亿元)"
$detail.Range("H1").Value = "仓位排名"

$detailRows = @(
    @("900010", "中信卓越成长两年持有期混合A", "133.02", "93.07", "2.60", "3.4585", 9),
    @("900090", "中信卓越成长两年持有期混合B", "86.95",  "93.07", "2.60", "2.2607", 9),
    @("009394", "银华同力精选混合",             "20.03",  "94.68", "5.36", "1.0736", 7),
    @("161838", "银华创业板两年定期开放混合",     "10.44",  "95.40", "5.22", "0.5450", 10),
    @("900100", "中信卓越成长两年持有期混合C", "6.91",   "93.07", "2.60", "0.1797", 9),
    @("005434", "鹏华睿投灵活配置混合",         "3.41",   "82.48", "2.96", "0.1009", 2),
    @("006048", "长城中证500指数增强A",         "4.45",   "92.64", "1.98", "0.0881", 10),
    @("007413", "长城中证500指数增强C",         "1.72",   "92.64", "1.98", "0.0341", 10),
    @("008778", "嘉实中证500指数增强A",         "0.93",   "93.42", "2.02", "0.0188", 5),
    @("008779", "嘉实中证500指数增强C",         "0.35",   "93.42", "2.02", "0.0071", 5)
)

$lastDetailRow = $detailRows.Length + 1

# Column B (fund code) and D-G (numeric-looking figures) are stored as text
# in the source data -- force text formatting on the whole block before
# assigning values so Excel doesn't auto-convert them to numbers (this also
# preserves leading zeros in fund codes).
$detail.Range("B2:B$lastDetailRow").NumberFormat = "@"
$detail.Range("D2:G$lastDetailRow").NumberFormat = "@"

for ($i = 0; $i -lt $detailRows.Length; $i++) {
    $r = $i + 2
    $row = $detailRows[$i]

    $detail.Range("A$r").Value = $i
    $detail.Range("B$r").Value = $row[0]
    $detail.Range("C$r").Value = $row[1]
    $detail.Range("D$r").Value = $row[2]
    $detail.Range("E$r").Value = $row[3]
    $detail.Range("F$r").Value = $row[4]
    $detail.Range("G$r").Value = $row[5]
    $detail.Range("H$r").Value = $row[6]
}

# Header row + index column share the bold / bordered / centered look used
# throughout the workbook.
$detailHeader = $detail.Range("B1:H1")
$detailHeader.Font.Bold = $true
$detailHeader.Borders.LineStyle = 1
$detailHeader.HorizontalAlignment = -4108
$detailHeader.VerticalAlignment = -4160

$detailIndex = $detail.Range("A2:A$lastDetailRow")
$detailIndex.Font.Bold = $true
$detailIndex.Borders.LineStyle = 1
$detailIndex.HorizontalAlignment = -4108
$detailIndex.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 2) Populate the new "总计" sheet: header + 2022-Q1 row + previous rows.
# ---------------------------------------------------------------------
$summary.Cells.Clear()

$summary.Range("B1").Value = "日期"
$summary.Range("C1").Value = "持有数量(只)"
$summary.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 10, 7.77),
    @("2021-Q4", 24, 13.61),
    @("2021-Q3", 17, 15.25),
    @("2021-Q2", 20, 17.09),
    @("2021-Q1", 19, 14.99),
    @("2020-Q4", 16, 5.95)
)

$lastSummaryRow = $summaryRows.Length + 1

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]

    $summary.Range("A$r").Value = $i
    $summary.Range("B$r").Value = $row[0]
    $summary.Range("C$r").Value = $row[1]
    $summary.Range("D$r").Value = $row[2]
}

$summaryHeader = $summary.Range("B1:D1")
$summaryHeader.Font.Bold = $true
$summaryHeader.Borders.LineStyle = 1
$summaryHeader.HorizontalAlignment = -4108
$summaryHeader.VerticalAlignment = -4160

$summaryIndex = $summary.Range("A2:A$lastSummaryRow")
$summaryIndex.Font.Bold = $true
$summaryIndex.Borders.LineStyle = 1
$summaryIndex.HorizontalAlignment = -4108
$summaryIndex.VerticalAlignment = -4160
